$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "V12_n_criminal_record"
$ws.Range("B25").Value = "V22_main_crime_comission_date"
$ws.Range("B31").Value = "V27_program_duration_cat"
$ws.Range("B32").Value = "V28_days_from_crime_to_program"
$ws.Range("B33").Value = "V29_program_duration"
$ws.Range("B34").Value = "V30_program_start"
$ws.Range("B35").Value = "V31_program_end"
$ws.Range("B126").Value = "V115_RECID2015_recid"
$ws.Range("B144").Value = "V132_RECID2013_recid"
